# Commit: "Restored from revision ... TEST Author: admin. Type: SAVE."
# Diff shows cell C10 on sheet "Rules" changing from 18 to 1
# (style/numFmt stays the same - s="20" is unchanged).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 1
